$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.824.83'
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').Value = '1.560.01'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '205.45'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.45%  '
$ws.Range('E6').Value = '  -1.09%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '21.56'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -3.09%  '
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('E10').Value = '  -0.83%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0863'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.06%  '
$ws.Range('D12').Value = '1.783.83'
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').Value = '1.562.77'
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('E14').Value = '  -1.38%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.513'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.65%  '
$ws.Range('D16').Value = '26.832.58'
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '61.28'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -2.71%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '214.96'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.33'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.29%  '
$ws.Range('D20').Value = '0.0₃0682'
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.12'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.14'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.10%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.01'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.07%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '153.73'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.31%  '
$ws.Range('E26').Value = '  -0.73%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '14.99'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.38%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('E29').Value = '  -1.17%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0465'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.82%  '
$ws.Range('E31').Value = '  -3.45%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.19'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.67%  '
$ws.Range('D33').Value = '1.377.50'
$ws.Range('E33').Value = '  -1.13%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.91'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.23%  '
$ws.Range('E35').Value = '  -2.42%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.919'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.22%  '
$ws.Range('E38').Value = '  -0.91%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.525'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.80%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.808'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.56%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.991'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('E43').Value = '  +4.40%  '
$ws.Range('E44').Value = '  -0.54%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.18'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.89%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '63.50'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').Value = '1.696.36'
$ws.Range('E47').Value = '  -0.08%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '86.44'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.05%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0510'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.68%  '
$ws.Range('D50').Value = '0.0₇0967'
$ws.Range('E50').Value = '  -2.15%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0951'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.57%  '
